$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Give the two brand-new header cells the same "header" formatting (bold font + fill)
#    that the existing header row already uses, by format-painting from D1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2. Switch the numeric-looking columns to Text format so the values aren't reinterpreted as numbers.
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# 3. Populate the new data (TaxRate / TotalAmountIncludingTax columns).
$ws.Range("D2").Value = "1440"
$ws.Range("E1").Value = "TaxRate"
$ws.Range("F1").Value = "TotalAmountIncludingTax"
$ws.Range("E2").Value = "0.0225"
$ws.Range("F2").Value = "1472"

# 4. Add a thin box border around every populated cell of the table.
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("E1").Borders.LineStyle = 1
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("B2").Borders.LineStyle = 1
$ws.Range("D2").Borders.LineStyle = 1
$ws.Range("E2").Borders.LineStyle = 1
$ws.Range("F2").Borders.LineStyle = 1

# 5. Column-level Text format + widths for the new/resized columns.
$ws.Columns("E").NumberFormat = "@"
$ws.Columns("F").NumberFormat = "@"
$ws.Columns("E").ColumnWidth = 17.67
$ws.Columns("F").ColumnWidth = 23.33

# 6. Move the selection like the author's last save.
$ws.Range("F2").Select() | Out-Null
